# Actualización automática 2025-10-01 08:30:09
#
# Monthly rollover:
#  - "VENTAS POR GRUPO" held the current (now-oldest) month's sales broken
#    down by product group; that month has rolled off the report so every
#    figure resets to 0 (and the "N de 58" completion counters at the
#    bottom follow suit).
#  - "VENTA MENSUAL" held four trailing months (junio..septiembre) per
#    client; everything shifts one column to the left (junio drops off,
#    julio -> junio's old slot, etc.) and a brand-new, still-empty
#    "octubre" column appears on the right, taking over column F.
#
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # VENTAS POR GRUPO
$ws2 = $wb.Worksheets.Item(2)   # VENTA MENSUAL

# ---------------------------------------------------------------------
# 1) VENTAS POR GRUPO: zero out every client/product figure (rows 2-59,
#    columns C-R), then reset the "N de 58" tally row (row 60) to zero.
# ---------------------------------------------------------------------
$firstDataRow = 2
$lastDataRow  = 59
$firstCol     = 3   # C
$lastCol      = 18  # R

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws1.Cells.Item($r, $c).Value = 0
    }
}

$totalRow = 60
for ($c = $firstCol; $c -le $lastCol; $c++) {
    $ws1.Cells.Item($totalRow, $c).Value = "0 de 58"
}

# ---------------------------------------------------------------------
# 2) VENTA MENSUAL: roll the four month columns (C:F) one slot to the
#    left across the header and every data/total row, then blank out
#    the newly-revealed rightmost month.
# ---------------------------------------------------------------------
$monthFirstCol = 3   # C
$monthLastCol  = 6   # F
$lastRow       = 60  # includes the totals row

$newMonthNames = @("julio", "agosto", "septiembre", "octubre")

# Header row (row 1) holds the month names as text.
for ($c = $monthFirstCol; $c -le $monthLastCol; $c++) {
    $ws2.Cells.Item(1, $c).Value = $newMonthNames[$c - $monthFirstCol]
}

# Data + totals rows (2..60): shift values left by one column, then
# zero the freshly-exposed rightmost column.
for ($r = 2; $r -le $lastRow; $r++) {
    $shifted = @()
    for ($c = $monthFirstCol + 1; $c -le $monthLastCol; $c++) {
        $shifted += $ws2.Cells.Item($r, $c).Value()
    }
    for ($i = 0; $i -lt $shifted.Count; $i++) {
        $ws2.Cells.Item($r, $monthFirstCol + $i).Value = $shifted[$i]
    }
    $ws2.Cells.Item($r, $monthLastCol).Value = 0
}

# Column widths (C:F) shift left the same way the data did; the new
# rightmost ("octubre") column reuses the width the old "julio" column
# (D) had, matching the target layout exactly.
$oldWidths = @()
for ($c = $monthFirstCol; $c -le $monthLastCol; $c++) {
    $oldWidths += $ws2.Columns.Item($c).ColumnWidth()
}
for ($c = $monthFirstCol; $c -lt $monthLastCol; $c++) {
    $ws2.Columns.Item($c).ColumnWidth = $oldWidths[$c - $monthFirstCol + 1]
}
$ws2.Columns.Item($monthLastCol).ColumnWidth = $oldWidths[1]
